# Generate Report for Archive
#
# Updates the localization-status report:
#   1. The "Status" value "Ready for handoff" becomes "In Translation"
#      (shared by the Overview sheet's zh-cn/de-de status columns and the
#      per-language "Status" column on each language sheet).
#   2. The two now-narrower "Status" columns are re-sized to match the
#      shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- 1. Update every cell currently showing the old status text ------------

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- 2. Shrink the Status columns to fit the new, shorter text -------------
# (17.2159881591797 -> 13.4101845877511 raw column-width units)

$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
